# "Add a little Purpose"
#
# 1. Drop the _GoBack bookmark that currently sits after the title
#    "PURPOSE AND OBJECTIVES".
# 2. Append a sentence to the (previously empty) "Purpose: " paragraph.
# 3. Re-create the _GoBack bookmark right after the new sentence (this is
#    where Word leaves it after the last edit).
# 4. Leave a few trailing spaces after the bookmark, as in the authored
#    edit.

$d = $word.ActiveDocument

# --- Step 1: remove the old _GoBack bookmark -----------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# --- Step 2: locate the "Purpose: " paragraph and append the new text ----
# (the one whose text is just the "Purpose:" heading run followed by ": ")
$purposePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^\s*Purpose:\s*$") {
        $purposePara = $p
        break
    }
}

$insertionPoint = $purposePara.Range
$insertionPoint.End = $insertionPoint.End - 1   # exclude the paragraph mark

$sentence = "Through this project, we can learn how to build database systems with relative sophisticated relationships. In the process, we will also be learning several techniques, including Java Servlet, JDBC, etc."
$insertionPoint.InsertAfter($sentence)

# Remember where the sentence ends - that's where the bookmark belongs.
$bookmarkPos = $insertionPoint.End

# --- Step 4: trailing spaces after the (future) bookmark -----------------
$insertionPoint.InsertAfter("   ")

# --- Step 3: re-create the _GoBack bookmark at the remembered position ---
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
